$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Split the "Profession: ..." paragraph into two paragraphs; the
#    first ("Profession: ") becomes a Heading1, the second keeps the
#    rest of the sentence as a normal paragraph.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Profession: Student aiming for a career in software engineering",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Profession: ^pStudent aiming for a career in software engineering", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Profession:*") {
        $p.Style = "Heading1"
        break
    }
}

# ---------------------------------------------------------------------
# 2. Rewrite the "About Me" paragraph text. (Assigned directly via
#    Range.Text rather than Find/Replace so the straight apostrophe in
#    "Bachelor's" is not auto-corrected into a curly quote.)
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "I am a student finishing my DEC*") {
        $p.Range.Text = "I am a student completing my DEC in computer science and mathematics. I have a strong interest in technology and computing. I aim to continue my development by enrolling in a Bachelor's degree in software engineering at E.T.S."
        break
    }
}

# ---------------------------------------------------------------------
# 3. Skills section: drop the bullet-list style and replace the seven
#    bullets with a reordered / updated set of plain paragraphs.
# ---------------------------------------------------------------------
$skillsHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Skills*" -and $p.Style.NameLocal -eq "Heading 1") {
        $skillsHeading = $p
        break
    }
}

$skillsStart = $skillsHeading.Index + 1
$skillsEnd = $skillsStart
while ($d.Paragraphs.Item($skillsEnd).Style.NameLocal -eq "List Bullet") {
    $skillsEnd = $skillsEnd + 1
}
$skillsEnd = $skillsEnd - 1

for ($i = $skillsEnd; $i -ge $skillsStart; $i--) {
    $d.Paragraphs.Item($i).Range.Delete()
}

$newSkills = @(
    "- Strong analytical and problem-solving abilities",
    "- Web design and development",
    "- Java programming skills",
    "- Requirements gathering and functional specification",
    "- Team coordination and leadership",
    "- Communication and collaboration",
    "- Python development and automation"
)
$insertAfter = $skillsHeading
foreach ($line in $newSkills) {
    $insertAfter.Range.InsertParagraphAfter()
    $insertAfter = $insertAfter.Next()
    $insertAfter.Range.Text = $line
    $insertAfter.Style = "Normal"
}

# ---------------------------------------------------------------------
# 4. Experience section: strip bullet-list style + update wording.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "List Bullet" -and $p.Range.Text -like "*Cook at McDonald*") {
        $p.Style = "Normal"
    }
    elseif ($p.Style.NameLocal -eq "List Bullet" -and $p.Range.Text -like "*Dockworker at Saint-Mathieux-De-Rioux*") {
        $p.Style = "Normal"
    }
    elseif ($p.Style.NameLocal -eq "List Bullet" -and $p.Range.Text -like "*Clerk at*Vaudreuil*") {
        $p.Style = "Normal"
    }
    elseif ($p.Style.NameLocal -eq "List Bullet" -and $p.Range.Text -like "*Cook at Premi*re Moisson*") {
        $p.Style = "Normal"
    }
}

$d.Content.Find.Execute(
    "Fast food preparation, adherence to hygiene standards, efficiency during peak hours, adaptability, and teamwork.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Quick food preparation, hygiene standards compliance, efficiency during peak hours, adaptability, and teamwork.", 2) | Out-Null

$d.Content.Find.Execute(
    "Dockworker at Saint-Mathieux-De-Rioux: Effective management of port operations, user assistance, maintenance of space, and application of local regulations.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Dock attendant at Saint-Mathieux-De-Rioux: Efficient management of port operations, assistance to users, maintenance of space, and enforcement of local regulations.", 2) | Out-Null

$d.Content.Find.Execute(
    "Precise preparation of orders in warehouse, inventory management, and optimization of storage space.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Precise preparation of warehouse orders, inventory management, and space optimization.", 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Education section: strip bullet-list style + append date ranges.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "List Bullet" -and (
            $p.Range.Text -like "*Elementary school*" -or
            $p.Range.Text -like "*High school*" -or
            $p.Range.Text -like "*CEGEP*")) {
        $p.Style = "Normal"
    }
}

$d.Content.Find.Execute(
    "Elementary school: Virginie-Roy",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Elementary School: Virginie-Roy 2021-2024", 2) | Out-Null

$d.Content.Find.Execute(
    "High school: Collège Notre-Dame",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "High School: Collège Notre-Dame 2016-2021", 2) | Out-Null

$d.Content.Find.Execute(
    "CEGEP: Cégep Gérald Godin",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CEGEP: Cégep Gérald Godin 2009-2016", 2) | Out-Null

# ---------------------------------------------------------------------
# 6. Programming Languages & Projects: strip bullet-list style + update
#    wording.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "List Bullet" -and (
            $p.Range.Text -like "*Website designer*" -or
            $p.Range.Text -like "*Java skills*" -or
            $p.Range.Text -like "*Exploring Python*")) {
        $p.Style = "Normal"
    }
}

$d.Content.Find.Execute(
    "Website designer: Creation of a professional CV site, implementation of attractive designs and optimization for better online visibility.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Website Designer: Created a professional CV site, implemented attractive designs, and optimized for better online visibility.", 2) | Out-Null

$d.Content.Find.Execute(
    "Java skills: Proficiency in Java and use of JavaFX interfaces and Scene Builder for the development of interactive applications.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Java Skills: Proficient in Java and used JavaFX and Scene Builder interfaces for developing interactive applications.", 2) | Out-Null

$d.Content.Find.Execute(
    "Exploring Python: Development of an operational Discord bot to perform essential tasks, implementing programming and automation skills.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Exploring Python: Developed a functional Discord bot to execute essential tasks, implementing programming and automation skills.", 2) | Out-Null

# ---------------------------------------------------------------------
# 7. Contact Details: merge the Email/Phone lines into a single
#    paragraph separated by a line break, and drop the Social Media
#    paragraph entirely.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Social Media:*") {
        $p.Range.Delete()
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Phone:*") {
        $p.Range.Delete()
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Email:*") {
        $p.Range.InsertAfter("`vPhone: 514-451-6262")
        break
    }
}
